$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1083
$ws.Range("F4").Value = 1138
$ws.Range("F5").Value = 397
$ws.Range("F6").Value = 153
$ws.Range("F7").Value = 514
$ws.Range("F8").Value = 254
$ws.Range("F9").Value = 47
$ws.Range("F10").Value = 1216
$ws.Range("F11").Value = 27585
$ws.Range("F12").Value = 3095
$ws.Range("F13").Value = 26
$ws.Range("F14").Value = 222
$ws.Range("F15").Value = 428
$ws.Range("F17").Value = 279
$ws.Range("F18").Value = 547
$ws.Range("F19").Value = 250
$ws.Range("F20").Value = 219
$ws.Range("F21").Value = 327
$ws.Range("F23").Value = 635
$ws.Range("F24").Value = 170
$ws.Range("F25").Value = 77
$ws.Range("F26").Value = 468
$ws.Range("F27").Value = 50
$ws.Range("F28").Value = 31
$ws.Range("F29").Value = 562
$ws.Range("F31").Value = 26

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 197
$ws.Range("F4").Value = 14
$ws.Range("F6").Value = 352
$ws.Range("F7").Value = 669
$ws.Range("F8").Value = 71
$ws.Range("F9").Value = 261
$ws.Range("F10").Value = 4207
$ws.Range("F12").Value = 176
$ws.Range("F15").Value = 33
$ws.Range("F17").Value = 33
$ws.Range("F18").Value = 44
$ws.Range("F21").Value = 4200

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 231
$ws.Range("F4").Value = 1093
$ws.Range("F5").Value = 284

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 231
$ws.Range("F4").Value = 1093
$ws.Range("F6").Value = 197
$ws.Range("F7").Value = 14
$ws.Range("F9").Value = 352
$ws.Range("F10").Value = 284
$ws.Range("F11").Value = 669
$ws.Range("F12").Value = 1083
$ws.Range("F13").Value = 1138
$ws.Range("F14").Value = 153
$ws.Range("F15").Value = 514
$ws.Range("F16").Value = 254
$ws.Range("F17").Value = 47
$ws.Range("F18").Value = 1216
$ws.Range("F19").Value = 27586
$ws.Range("F20").Value = 71
$ws.Range("F21").Value = 261
$ws.Range("F23").Value = 176
$ws.Range("F26").Value = 3095
$ws.Range("F27").Value = 222
$ws.Range("F28").Value = 33
$ws.Range("F29").Value = 33
$ws.Range("F30").Value = 428
$ws.Range("F32").Value = 33
$ws.Range("F33").Value = 279
$ws.Range("F34").Value = 547
$ws.Range("F35").Value = 250
$ws.Range("F36").Value = 327
$ws.Range("F38").Value = 635
$ws.Range("F39").Value = 44
$ws.Range("F40").Value = 170
$ws.Range("F41").Value = 77
$ws.Range("F44").Value = 50
$ws.Range("F45").Value = 31
$ws.Range("F46").Value = 562
$ws.Range("F48").Value = 4200
$ws.Range("F49").Value = 26

Write-Output "Updated F-column values across all sheets."
